$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the slightly-adjusted timestamp on row 18 (A18)
$ws.Range("A18").Value = 45818.3937591088

# Append new rows 19-37 with the same product/weight/price pattern
$ws.Range("A19").Value = 45833.41288174769
$ws.Range("A19").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("B19").Value = "CREATINA MONOHIDRATO EN POLVO"
$ws.Range("C19").Value = "1Kg"
$ws.Range("D19").Value = "15,41€"

$ws.Range("A20").Value = 45833.41722056713
$ws.Range("A20").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("B20").Value = "CREATINA MONOHIDRATO EN POLVO"
$ws.Range("C20").Value = "1Kg"
$ws.Range("D20").Value = "15,41€"

$ws.Range("A21").Value = 45833.41746144676
$ws.Range("A21").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("B21").Value = "CREATINA MONOHIDRATO EN POLVO"
$ws.Range("C21").Value = "1Kg"
$ws.Range("D21").Value = "15,41€"

$ws.Range("A22").Value = 45833.41859150463
$ws.Range("A22").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("B22").Value = "CREATINA MONOHIDRATO EN POLVO"
$ws.Range("C22").Value = "1Kg"
$ws.Range("D22").Value = "15,41€"

$ws.Range("A23").Value = 45833.42342791666
$ws.Range("A23").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("B23").Value = "CREATINA MONOHIDRATO EN POLVO"
$ws.Range("C23").Value = "1Kg"
$ws.Range("D23").Value = "15,41€"

$ws.Range("A24").Value = 45833.42990267361
$ws.Range("A24").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("B24").Value = "CREATINA MONOHIDRATO EN POLVO"
$ws.Range("C24").Value = "1Kg"
$ws.Range("D24").Value = "15,41€"

$ws.Range("A25").Value = 45833.43051061343
$ws.Range("A25").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("B25").Value = "CREATINA MONOHIDRATO EN POLVO"
$ws.Range("C25").Value = "1Kg"
$ws.Range("D25").Value = "15,41€"

$ws.Range("A26").Value = 45833.43422083333
$ws.Range("A26").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("B26").Value = "CREATINA MONOHIDRATO EN POLVO"
$ws.Range("C26").Value = "1Kg"
$ws.Range("D26").Value = "15,41€"

$ws.Range("A27").Value = 45833.4369330787
$ws.Range("A27").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("B27").Value = "CREATINA MONOHIDRATO EN POLVO"
$ws.Range("C27").Value = "1Kg"
$ws.Range("D27").Value = "15,41€"

$ws.Range("A28").Value = 45833.43815046296
$ws.Range("A28").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("B28").Value = "CREATINA MONOHIDRATO EN POLVO"
$ws.Range("C28").Value = "1Kg"
$ws.Range("D28").Value = "15,41€"

$ws.Range("A29").Value = 45833.43892174769
$ws.Range("A29").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("B29").Value = "CREATINA MONOHIDRATO EN POLVO"
$ws.Range("C29").Value = "1Kg"
$ws.Range("D29").Value = "15,41€"

$ws.Range("A30").Value = 45833.43977028935
$ws.Range("A30").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("B30").Value = "CREATINA MONOHIDRATO EN POLVO"
$ws.Range("C30").Value = "1Kg"
$ws.Range("D30").Value = "15,41€"

$ws.Range("A31").Value = 45833.44140856482
$ws.Range("A31").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("B31").Value = "CREATINA MONOHIDRATO EN POLVO"
$ws.Range("C31").Value = "1Kg"
$ws.Range("D31").Value = "15,41€"

$ws.Range("A32").Value = 45833.44283475694
$ws.Range("A32").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("B32").Value = "CREATINA MONOHIDRATO EN POLVO"
$ws.Range("C32").Value = "1Kg"
$ws.Range("D32").Value = "15,41€"

$ws.Range("A33").Value = 45833.44410322917
$ws.Range("A33").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("B33").Value = "CREATINA MONOHIDRATO EN POLVO"
$ws.Range("C33").Value = "1Kg"
$ws.Range("D33").Value = "15,41€"

$ws.Range("A34").Value = 45833.44530853009
$ws.Range("A34").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("B34").Value = "CREATINA MONOHIDRATO EN POLVO"
$ws.Range("C34").Value = "1Kg"
$ws.Range("D34").Value = "15,41€"

$ws.Range("A35").Value = 45833.45332490741
$ws.Range("A35").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("B35").Value = "CREATINA MONOHIDRATO EN POLVO"
$ws.Range("C35").Value = "1Kg"
$ws.Range("D35").Value = "15,41€"

$ws.Range("A36").Value = 45833.45424282407
$ws.Range("A36").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("B36").Value = "CREATINA MONOHIDRATO EN POLVO"
$ws.Range("C36").Value = "1Kg"
$ws.Range("D36").Value = "15,41€"

$ws.Range("A37").Value = 45833.45891689683
$ws.Range("A37").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("B37").Value = "CREATINA MONOHIDRATO EN POLVO"
$ws.Range("C37").Value = "1Kg"
$ws.Range("D37").Value = "15,41€"
